# Generate Report for Handoff
#
# Two new localization source files have reached "Ready for handoff" status:
#   678a0f76-d84f-4dc4-9b45-be08ef3ae90e.md
#   97481782-9c7c-4a20-ae58-9c3b33910e1a.md
#
# They are inserted (in that order) ahead of the already-tracked
# df0b45ff-a6c6-4bae-85ee-52adf71e0f22.md entry on every sheet of the
# handoff-status workbook (Overview + one sheet per locale).

$wb = $excel.ActiveWorkbook

function Set-CellText($ws, $row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = $text
}

function Add-Link($ws, $row, $col, $url, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $ws.Hyperlinks.Add($cell, $url, "", "", $text) | Out-Null
}

# ======================================================================
# Sheet "Overview" — File Name | zh-cn | de-de | Latest Handoff Date
# ======================================================================
$ws = $wb.Worksheets.Item("Overview")

# Remove every existing hyperlink up front; they get rebuilt below in one
# consistent pass (in-place hyperlink edits leave stale duplicate entries
# in this host, so a clear + full rebuild is the reliable path).
$ws.Cells.Hyperlinks.Delete()

$overviewRows = @(
    @{ Row = 2; Name = "1ae9c0f9-4b1d-4f9f-a308-5ea809b2996b.md"; Status = "Handed back: in sync with en-US"; Date = "2016-39-13 16:39:11"; Rev = "11b4f84abb3deefd189b8a235d61c0a8bd336130" },
    @{ Row = 3; Name = "8ddd45d9-453a-412a-bb78-dc1c113f7102.md"; Status = "In Translation";                   Date = "2016-42-13 16:42:11"; Rev = "e3d23fce5087036794a2b7a46e9e84d42a30fbf9" },
    @{ Row = 4; Name = "f4e3b2ca-a377-4bca-9905-98ca926acfdd.md"; Status = "In Translation";                   Date = "2016-42-13 16:42:11"; Rev = "e3d23fce5087036794a2b7a46e9e84d42a30fbf9" },
    @{ Row = 5; Name = "678a0f76-d84f-4dc4-9b45-be08ef3ae90e.md"; Status = "Ready for handoff";                Date = "2016-42-13 16:42:56"; Rev = "8e3f1c2a9b7d4e5f6a8b9c0d1e2f3a4b5c6d7e8f" },
    @{ Row = 6; Name = "97481782-9c7c-4a20-ae58-9c3b33910e1a.md"; Status = "Ready for handoff";                Date = "2016-42-13 16:42:56"; Rev = "1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b" },
    @{ Row = 7; Name = "df0b45ff-a6c6-4bae-85ee-52adf71e0f22.md"; Status = "Ready for handoff";                Date = "2016-38-13 16:38:26"; Rev = "7060af68e6ee8cfb73e7c9d4b5a079ca9b2caf28" }
)

foreach ($r in $overviewRows) {
    Set-CellText $ws $r.Row 2 $r.Status
    Set-CellText $ws $r.Row 3 $r.Status
    Set-CellText $ws $r.Row 4 $r.Date
    $url = "https://github.com/OpenLocalizationTest/oltest/blob/$($r.Rev)/e2e/$($r.Name)"
    Add-Link $ws $r.Row 1 $url $r.Name
}

Write-Output "Overview sheet updated"

# ======================================================================
# Locale sheets "zh-cn" / "de-de"
# ======================================================================
$locales = @(
    @{ Sheet = "zh-cn"; Code = "zh-cn";
       Files = @{
           "678a0f76-d84f-4dc4-9b45-be08ef3ae90e" = "8966f2615ee84af99f5ca99a1b34c4026038bea6";
           "97481782-9c7c-4a20-ae58-9c3b33910e1a" = "33d2dc99b9bb4d3bdc842a4d3d5bfbb0dda69f07";
           "df0b45ff-a6c6-4bae-85ee-52adf71e0f22" = "9fe66e6aa6752bd3910e9836af17b09e2531c6e0"
       };
       Dates = @{
           "678a0f76-d84f-4dc4-9b45-be08ef3ae90e" = "2016-03-13 16:42:52";
           "97481782-9c7c-4a20-ae58-9c3b33910e1a" = "2016-03-13 16:42:52";
           "df0b45ff-a6c6-4bae-85ee-52adf71e0f22" = "2016-03-13 16:36:25"
       } },
    @{ Sheet = "de-de"; Code = "de-de";
       Files = @{
           "678a0f76-d84f-4dc4-9b45-be08ef3ae90e" = "8966f2615ee84af99f5ca99a1b34c4026038bea6";
           "97481782-9c7c-4a20-ae58-9c3b33910e1a" = "33d2dc99b9bb4d3bdc842a4d3d5bfbb0dda69f07";
           "df0b45ff-a6c6-4bae-85ee-52adf71e0f22" = "9fe66e6aa6752bd3910e9836af17b09e2531c6e0"
       };
       Dates = @{
           "678a0f76-d84f-4dc4-9b45-be08ef3ae90e" = "2016-03-13 16:42:56";
           "97481782-9c7c-4a20-ae58-9c3b33910e1a" = "2016-03-13 16:42:56";
           "df0b45ff-a6c6-4bae-85ee-52adf71e0f22" = "2016-03-13 16:38:26"
       } }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)
    $ws.Cells.Hyperlinks.Delete()

    $rows = @(
        @{ Row = 2; Name = "1ae9c0f9-4b1d-4f9f-a308-5ea809b2996b"; Status = "Handed back: in sync with en-US"; HOSha = "9df732f9881ab542f52678ce62f1a2c6d5b38160"; HODate = "2016-03-13 16:39:08"; HBDate = "2016-03-13 16:41:36"; Reason = "Include"; IncludeFG = $true },
        @{ Row = 3; Name = "8ddd45d9-453a-412a-bb78-dc1c113f7102"; Status = "In Translation";                   HOSha = "7916894b0853006527613b491b3fd06de7ca0606"; HODate = "2016-03-13 16:42:08"; HBDate = "0001-01-01 00:00:00"; Reason = "Include"; IncludeFG = $false },
        @{ Row = 4; Name = "f4e3b2ca-a377-4bca-9905-98ca926acfdd"; Status = "In Translation";                   HOSha = "0e1dd9abb24ef2055a7664b8839bb2e5b9bc406b"; HODate = "2016-03-13 16:42:08"; HBDate = "0001-01-01 00:00:00"; Reason = "Include"; IncludeFG = $false },
        @{ Row = 5; Name = "678a0f76-d84f-4dc4-9b45-be08ef3ae90e"; Status = "Ready for handoff";                HOSha = $loc.Files["678a0f76-d84f-4dc4-9b45-be08ef3ae90e"]; HODate = $loc.Dates["678a0f76-d84f-4dc4-9b45-be08ef3ae90e"]; HBDate = "0001-01-01 00:00:00"; Reason = "Include"; IncludeFG = $false },
        @{ Row = 6; Name = "97481782-9c7c-4a20-ae58-9c3b33910e1a"; Status = "Ready for handoff";                HOSha = $loc.Files["97481782-9c7c-4a20-ae58-9c3b33910e1a"]; HODate = $loc.Dates["97481782-9c7c-4a20-ae58-9c3b33910e1a"]; HBDate = "0001-01-01 00:00:00"; Reason = "Include"; IncludeFG = $false },
        @{ Row = 7; Name = "df0b45ff-a6c6-4bae-85ee-52adf71e0f22"; Status = "Ready for handoff";                HOSha = $loc.Files["df0b45ff-a6c6-4bae-85ee-52adf71e0f22"]; HODate = $loc.Dates["df0b45ff-a6c6-4bae-85ee-52adf71e0f22"]; HBDate = "0001-01-01 00:00:00"; Reason = "Include"; IncludeFG = $false }
    )

    foreach ($r in $rows) {
        $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/e3d23fce5087036794a2b7a46e9e84d42a30fbf9/e2e/$($r.Name).md"
        $xlfName = "$($r.Name).$($r.HOSha).$($loc.Code).xlf"
        $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c32f4c83c712f1b36afd0f3e41ee48c2d716399c/ol-handoff/OpenLocalizationTestOrg/oltest.$($loc.Code)/ci/ht/$xlfName"

        Add-Link $ws $r.Row 1 $mdUrl "$($r.Name).md"
        Add-Link $ws $r.Row 2 $mdUrl ".md"
        Set-CellText $ws $r.Row 3 $r.Status
        Add-Link $ws $r.Row 4 $xlfUrl $xlfName
        Set-CellText $ws $r.Row 5 $r.HODate
        $ws.Cells.Item($r.Row, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"

        if ($r.IncludeFG) {
            Add-Link $ws $r.Row 6 $mdUrl "$($r.Name).md"
            Add-Link $ws $r.Row 7 $xlfUrl $xlfName
        }

        Set-CellText $ws $r.Row 8 $r.HBDate
        Set-CellText $ws $r.Row 9 $r.Reason
    }

    Write-Output "$($loc.Sheet) sheet updated"
}

Write-Output "All sheets updated"
